# Updated symbol list on Tue Dec 13 16:43:54 UTC 2022 with GitHub Actions
# Applies the price/label refresh to the "Price" (D) and "Volume(1h)" (E)
# columns on the active sheet, preserving each cell's original Text type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Value
    )
    $rng = $ws.Range($Address)
    # Force the cell to Text format so numeric-looking strings (e.g. "274.51")
    # are NOT auto-converted to a number, then restore the default "Normal"
    # style so no stray number-format/style is left behind on save.
    $rng.NumberFormat = "@"
    $rng.Value = $Value
    $rng.Style = "Normal"
}

Set-TextValue "D2"  "274.51"
Set-TextValue "D3"  "22.94"
Set-TextValue "D4"  "6.307"
Set-TextValue "D5"  "0.06236"
Set-TextValue "D6"  "3.644"
Set-TextValue "D7"  "6.675"
Set-TextValue "D8"  "1.398"
Set-TextValue "D9"  "0.8330"
Set-TextValue "D11" "0.1604"
Set-TextValue "D12" "0.08326"
Set-TextValue "D13" "0.03541"
Set-TextValue "D15" "4.062"
Set-TextValue "D16" "0.09305"
Set-TextValue "D17" "0.001654"
Set-TextValue "D18" "0.04730"
Set-TextValue "D19" "0.006291"
Set-TextValue "D20" "0.005710"
Set-TextValue "D23" "3.732"
Set-TextValue "D25" "0.3361"
Set-TextValue "D26" "0.1260"
Set-TextValue "D27" "0.0002702"
Set-TextValue "D40" "0.04739"
Set-TextValue "D41" "0.006921"
Set-TextValue "D42" "0.003900"
Set-TextValue "E42" "41CEJICEJIWorstin24h"
Set-TextValue "D43" "0.1167"
Set-TextValue "D45" "0.00006094"
Set-TextValue "D46" "0.0009895"
Set-TextValue "D48" "0.7817"
Set-TextValue "E48" "47CoinbaseStockTokenCOIN"
Set-TextValue "D49" "0.002358"
Set-TextValue "D50" "0.00002400"
